# Applies the "Improve data source parsing and dashboard visuals" edit:
#  - Removes the "Petroleo Brasileiro Petrobras SA ADR" and "Vale SA ADR" rows
#    (their marketwatch.com data sources were dropped).
#  - Flips the sign of the base values for the currency-pair / index rows
#    that now sit above the ETF block (EUR/BRL ... Indice Dolar Futuros).
#  - Repoints the USD/KRW data source from marketwatch.com to investing.com.
#  - Rebuilds the hyperlink list so it lines up with the now-shifted rows.
#  - Updates the active-cell selection to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the two discontinued rows (Petrobras ADR, Vale ADR). Everything
#    below shifts up by two rows automatically.
# ---------------------------------------------------------------------
$ws.Rows("24:25").Delete()

# ---------------------------------------------------------------------
# 2. Negate the base values for the currency/index rows now in rows 2-10
#    (EUR/BRL, USD/CNY, USD/KRW, USD/AUD, USD/NZD, USD/NOK, USD/MXN,
#    S&P 500 VIX Futuros, Indice Dolar Futuros).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = -$cell.Value()
}

# ---------------------------------------------------------------------
# 3. Switch the USD/KRW data source URL to investing.com.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "https://br.investing.com/indices/f-usdkrw"

# ---------------------------------------------------------------------
# 4. Rebuild the hyperlinks: clear the stale set (refs now point at the
#    wrong rows after the delete) and re-add them against the current
#    layout. C4 (USD/KRW) and C24 (Global Dow) stay plain text, matching
#    the original workbook.
# ---------------------------------------------------------------------
$ws.Range("A1:C31").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "https://br.investing.com/currencies/eur-brl")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://br.investing.com/currencies/usd-cny")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://br.investing.com/currencies/usd-aud")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://br.investing.com/currencies/usd-nzd")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://br.investing.com/currencies/usd-nok")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://br.investing.com/currencies/usd-mxn")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://br.investing.com/indices/us-spx-vix-futures")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://br.investing.com/currencies/us-dollar-index")
$ws.Hyperlinks.Add($ws.Range("C11"), "https://br.investing.com/indices/shanghai-composite")
$ws.Hyperlinks.Add($ws.Range("C12"), "https://br.investing.com/indices/szse-component")
$ws.Hyperlinks.Add($ws.Range("C13"), "https://br.investing.com/indices/dj-shanghai")
$ws.Hyperlinks.Add($ws.Range("C14"), "https://br.investing.com/indices/china-a50")
$ws.Hyperlinks.Add($ws.Range("C15"), "https://br.investing.com/indices/hong-kong-40-futures")
$ws.Hyperlinks.Add($ws.Range("C16"), "https://br.investing.com/indices/sensex")
$ws.Hyperlinks.Add($ws.Range("C17"), "https://br.investing.com/etfs/ishares-phlx-sox-semiconductor")
$ws.Hyperlinks.Add($ws.Range("C18"), "https://br.investing.com/etfs/ishares-msci-emg-markets")
$ws.Hyperlinks.Add($ws.Range("C19"), "https://br.investing.com/etfs/spdr-s-p-metals---mining")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://br.investing.com/etfs/spdr-energy-select-sector-fund")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://br.investing.com/etfs/spdr---consumer-staples")
$ws.Hyperlinks.Add($ws.Range("C22"), "https://br.investing.com/etfs/financial-select-sector-spdr-fund")
$ws.Hyperlinks.Add($ws.Range("C23"), "https://br.investing.com/etfs/ishares-brazil-index")
$ws.Hyperlinks.Add($ws.Range("C25"), "https://br.investing.com/indices/us-30-futures")
$ws.Hyperlinks.Add($ws.Range("C26"), "https://br.investing.com/commodities/us-soybeans")
$ws.Hyperlinks.Add($ws.Range("C27"), "https://br.investing.com/indices/oslo-all-share")
$ws.Hyperlinks.Add($ws.Range("C28"), "https://br.investing.com/commodities/crude-oil")
$ws.Hyperlinks.Add($ws.Range("C29"), "https://br.investing.com/commodities/copper")
$ws.Hyperlinks.Add($ws.Range("C30"), "https://br.investing.com/commodities/gold")

# ---------------------------------------------------------------------
# 5. Match the saved selection/active cell.
# ---------------------------------------------------------------------
$ws.Range("N12").Select()
